$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("testcase_1")
$ws2 = $wb.Worksheets.Item("testcase_2")

# --- testcase_1: insert a new "reference" column before column A ---
$ws1.Activate()
$ws1.Columns("A").Insert()

$ws1.Range("A1").Value = "reference"
$ws1.Range("A2").Value = "aaaaa"
$ws1.Range("A3").Value = "aaaaa"
$ws1.Range("A4").Value = "dfdfdef"
$ws1.Range("A2:A4").WrapText = $true

# --- testcase_2: just a selection/view change ---
$ws2.Activate()
$ws2.Range("R9").Select()

# restore testcase_1 as the active/selected sheet & cell
$ws1.Activate()
$ws1.Range("A4").Select()
